$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new rows of data for additional Fiat 500 entries (columns are
# still in the original B/C/D/E layout at this point).
$ws.Range("B8").Value = "Veicolo"
$ws.Range("C8").Value = "Fiat"
$ws.Range("D8").Value = 500
$ws.Range("E8").Value = 2015

$ws.Range("B9").Value = "Veicolo"
$ws.Range("C9").Value = "Fiat"
$ws.Range("D9").Value = 500
$ws.Range("E9").Value = 2020

$ws.Range("B10").Value = "Veicolo"
$ws.Range("C10").Value = "Fiat"
$ws.Range("D10").Value = 500
$ws.Range("E10").Value = 2022

$ws.Range("B11").Value = "Veicolo"
$ws.Range("C11").Value = "Fiat"
$ws.Range("D11").Value = "500L"
$ws.Range("E11").Value = 2022

$ws.Range("B12").Value = "Veicolo"
$ws.Range("C12").Value = "Fiat"
$ws.Range("D12").Value = "500 L"
$ws.Range("E12").Value = 2022

# Insert a new column before column C (shifts old C/D/E -> D/E/F) and
# fill in the new "VIN" column.
$ws.Columns("C:C").Insert()
$ws.Range("C1").Value = "VIN"
$ws.Range("C2").Value = "gvr"
$ws.Range("C4").Value = "grer"

# Restore the active selection as it was after editing.
$ws.Range("C4").Select() | Out-Null
